$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row J1:M1 - set values then copy the existing header formatting (style s=1)
$ws.Range("J1").Value = "Mean_T_anom_2002-2007"
$ws.Range("K1").Value = "Mean_T_anom_2007-2012"
$ws.Range("L1").Value = "Mean_T_anom_2012-2017"
$ws.Range("M1").Value = "Mean_T_anom_2017-2022"

$ws.Range("B1").Copy()
$ws.Range("J1:M1").PasteSpecial(-4122)

# Data rows J2:M15 (plain numeric cells, no special style - matches columns B:I in data rows)
$data = @(
    @(1.027534365653992, 0.5953571796417236, 0.9207971096038818, 1.407952666282654),
    @(0.1529367417097092, -0.2098842114210129, 0.1759693473577499, 0.476120263338089),
    @(0.3213476538658142, 0.2530762851238251, 0.4622536599636078, 0.8613607883453369),
    @(0.6876922845840454, 0.4942399561405182, 0.7095251679420471, 1.114517331123352),
    @(0.9773463010787964, 0.6907257437705994, 1.033843874931335, 1.462576985359192),
    @(0.8976433277130127, 0.5702357888221741, 0.7934548258781433, 1.276681780815125),
    @(0.9124077558517456, 0.5719060897827148, 0.8839142322540283, 1.381705760955811),
    @(0.8915310502052307, 0.4934071600437164, 0.9293677806854248, 1.325136542320251),
    @(0.1984381079673767, 0.1777880936861038, 0.4904373586177826, 0.5138602256774902),
    @(0.2728450894355774, 0.05762287601828575, 0.405795693397522, 0.7220320701599121),
    @(0.3603193163871765, 0.4083164036273956, 0.5908300876617432, 0.8542513251304626),
    @(0.7681716680526733, 0.5500335097312927, 0.9291219711303711, 1.46335780620575),
    @(0.6368283033370972, 0.9244036078453064, 1.339719653129578, 1.588810563087463),
    @(0.7681716680526733, 0.5500335097312927, 0.9291219711303711, 1.46335780620575)
)

$row = 2
foreach ($r in $data) {
    $ws.Cells.Item($row, 10).Value = $r[0]
    $ws.Cells.Item($row, 11).Value = $r[1]
    $ws.Cells.Item($row, 12).Value = $r[2]
    $ws.Cells.Item($row, 13).Value = $r[3]
    $row++
}
